$d = $word.ActiveDocument

# The "Requisitos" bullet list currently has three requirement lines,
# each its own run ending in a manual line break (<w:br/>).
# Replace all three with a single requirement line (also ending in a
# manual line break), collapsing them into one run.
$old = "LOM3036 -  Propriedades Mecânicas  (Requisito fraco)`vLOM3082 -  Cerâmica Física  (Requisito fraco)`vLOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)`v"
$new = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`v"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the expected 'Requisitos' text to replace."
}
